# 1. Login_eSign.xlsx — "Add files via upload"
#
# - Active tab moves from "Setting" to "Login".
# - "Login" sheet: the ADINS admin-esign login URL changes from
#   https://gdkwebserver.ad-ins.com/ to http://gdkwebsvr:8080/login
#   (used by rows 2-4, and by a brand-new row 5: admesign/password/ADINS/Admin Esign).
# - "Login" sheet C8 changes from "esign_uat" to "esign".
# - Selection on "Login" moves to C9; selection on "Setting" stays at G4
#   but "Setting" is no longer the selected tab.

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")
$wsSetting = $wb.Worksheets.Item("Setting")

# Seed row 5 with row 4's plain formatting (border style, no fill) before
# filling in values, so the new cells pick up the same style ("s=1") the
# rest of the data rows use.
$wsLogin.Range("A4:E4").Copy()
$wsLogin.Range("A5:E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 5 values (order matches the original authoring sequence so newly
# introduced shared strings land in the same table order as the source edit).
$wsLogin.Range("B5").Value = "admesign"
$wsLogin.Range("D5").Value = "ADINS"
$wsLogin.Range("E5").Value = "Admin Esign"
$wsLogin.Range("C5").Value = "password"

# Point the login URL cells (existing rows 2-4 plus new row 5) at the new
# admin-esign endpoint.
$wsLogin.Range("A2").Value = "http://gdkwebsvr:8080/login"
$wsLogin.Range("A3").Value = "http://gdkwebsvr:8080/login"
$wsLogin.Range("A4").Value = "http://gdkwebsvr:8080/login"
$wsLogin.Range("A5").Value = "http://gdkwebsvr:8080/login"

# Those URL cells use the same hyperlink-look style as Setting!D2 (underline
# + themed color, with the thin border); copy that formatting over rather
# than reassigning a named style so the existing style index is reused.
$wsSetting.Range("D2").Copy()
$wsLogin.Range("A2").PasteSpecial(-4122)
$wsLogin.Range("A3").PasteSpecial(-4122)
$wsLogin.Range("A4").PasteSpecial(-4122)
$wsLogin.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# esign_uat -> esign
$wsLogin.Range("C8").Value = "esign"

# Make "Login" the active sheet/tab, with C9 selected; "Setting" keeps its
# own G4 selection but is no longer the active tab.
$wsLogin.Activate()
$wsLogin.Range("C9").Select()
